$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.070051846520101
$ws.Range("D2").Value = 1.06626654332389
$ws.Range("E2").Value = 1.07338479428994
$ws.Range("F2").Value = 1.081096729699277
$ws.Range("I2").Value = 1.05443466022714
$ws.Range("J2").Value = 1.074983459293816
$ws.Range("K2").Value = 1.068977553544383
$ws.Range("L2").Value = 1.076076796518562
$ws.Range("M2").Value = 1.083768443738624
$ws.Range("N2").Value = 1.07651005910834
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.071326623303014
$ws.Range("D3").Value = 1.067048921584357
$ws.Range("E3").Value = 1.074463954362187
$ws.Range("F3").Value = 1.082122539378438
$ws.Range("I3").Value = 1.054773441197071
$ws.Range("J3").Value = 1.075914092385136
$ws.Range("K3").Value = 1.069575573768923
$ws.Range("L3").Value = 1.076972236134882
$ws.Range("M3").Value = 1.084612131521825
$ws.Range("N3").Value = 1.077442013805395
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.072151451856025
$ws.Range("D4").Value = 1.06755516541778
$ws.Range("E4").Value = 1.075162458477653
$ws.Range("F4").Value = 1.082786536972249
$ws.Range("I4").Value = 1.054991471227954
$ws.Range("J4").Value = 1.076515695605537
$ws.Range("K4").Value = 1.069961864986837
$ws.Range("L4").Value = 1.077551259262455
$ws.Range("M4").Value = 1.085157669968502
$ws.Range("N4").Value = 1.078044471371373
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.072498202769721
$ws.Range("D5").Value = 1.067767988441152
$ws.Range("E5").Value = 1.075456161942716
$ws.Range("F5").Value = 1.083065737066358
$ws.Range("I5").Value = 1.055082847863422
$ws.Range("J5").Value = 1.07676847209967
$ws.Range("K5").Value = 1.070124102059616
$ws.Range("L5").Value = 1.077794588897882
$ws.Range("M5").Value = 1.085386922884018
$ws.Range("N5").Value = 1.078297606837122
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.072556423418484
$ws.Range("D6").Value = 1.067803722248159
$ws.Range("E6").Value = 1.075505479138219
$ws.Range("F6").Value = 1.083112619257498
$ws.Range("I6").Value = 1.05509817380325
$ws.Range("J6").Value = 1.076810906352015
$ws.Range("K6").Value = 1.070151333013677
$ws.Range("L6").Value = 1.077835439660028
$ws.Range("M6").Value = 1.085425410108563
$ws.Range("N6").Value = 1.078340101350973
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.072156085183849
$ws.Range("D7").Value = 1.067558009177039
$ws.Range("E7").Value = 1.075166382751121
$ws.Range("F7").Value = 1.082790267438958
$ws.Range("I7").Value = 1.054992693319467
$ws.Range("J7").Value = 1.076519073756423
$ws.Range("K7").Value = 1.069964033434111
$ws.Range("L7").Value = 1.077554511004265
$ws.Range("M7").Value = 1.085160733615913
$ws.Range("N7").Value = 1.078047854319621
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.070482672168963
$ws.Range("D8").Value = 1.06653095312175
$ws.Range("E8").Value = 1.073749457355912
$ws.Range("F8").Value = 1.081443359479404
$ws.Range("I8").Value = 1.054549398123599
$ws.Range("J8").Value = 1.075298091526921
$ws.Range("K8").Value = 1.069179795642582
$ws.Range("L8").Value = 1.07637949495129
$ws.Range("M8").Value = 1.084053651763696
$ws.Range("N8").Value = 1.076825138155304
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.067533533492868
$ws.Range("D9").Value = 1.064721100156159
$ws.Range("E9").Value = 1.071254268430547
$ws.Range("F9").Value = 1.079071679210713
$ws.Range("I9").Value = 1.053759177009951
$ws.Range("J9").Value = 1.073142087205492
$ws.Range("K9").Value = 1.06779274716049
$ws.Range("L9").Value = 1.074305974233897
$ws.Range("M9").Value = 1.082099869674999
$ws.Range("N9").Value = 1.074666072060442
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.065567069489681
$ws.Range("D10").Value = 1.063514504458682
$ws.Range("E10").Value = 1.069591843161495
$ws.Range("F10").Value = 1.07749170194228
$ws.Range("I10").Value = 1.053226242468146
$ws.Range("J10").Value = 1.071701673770698
$ws.Range("K10").Value = 1.066864593873293
$ws.Range("L10").Value = 1.072921564465842
$ws.Range("M10").Value = 1.080795323657236
$ws.Range("N10").Value = 1.073223613073353
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.06471545289279
$ws.Range("D11").Value = 1.062992029172416
$ws.Range("E11").Value = 1.068872227587508
$ws.Range("F11").Value = 1.076807817955567
$ws.Range("I11").Value = 1.052994020087337
$ws.Range("J11").Value = 1.071077213410752
$ws.Range("K11").Value = 1.066461870787877
$ws.Range("L11").Value = 1.072321598481835
$ws.Range("M11").Value = 1.080229952507326
$ws.Range("N11").Value = 1.072598265908069
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.064399103340826
$ws.Range("D12").Value = 1.062797956809958
$ws.Range("E12").Value = 1.068604963044349
$ws.Range("F12").Value = 1.076553830795703
$ws.Range("I12").Value = 1.052907542854701
$ws.Range("J12").Value = 1.070845146637188
$ws.Range("K12").Value = 1.066312156959518
$ws.Range("L12").Value = 1.072098667035903
$ws.Range("M12").Value = 1.080019873458044
$ws.Range("N12").Value = 1.072365869573068
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.064466962394353
$ws.Range("D13").Value = 1.062839586072085
$ws.Range("E13").Value = 1.068662290734079
$ws.Range("F13").Value = 1.076608310205074
$ws.Range("I13").Value = 1.052926102454295
$ws.Range("J13").Value = 1.07089493094745
$ws.Range("K13").Value = 1.066344276725532
$ws.Range("L13").Value = 1.072146490106995
$ws.Range("M13").Value = 1.080064939535369
$ws.Range("N13").Value = 1.072415724582761
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.064689303787648
$ws.Range("D14").Value = 1.062975987127069
$ws.Range("E14").Value = 1.068850134756216
$ws.Range("F14").Value = 1.076786822530876
$ws.Range("I14").Value = 1.052986876332496
$ws.Range("J14").Value = 1.071058033040731
$ws.Range("K14").Value = 1.066449497938843
$ws.Range("L14").Value = 1.07230317248377
$ws.Range("M14").Value = 1.080212588842801
$ws.Range("N14").Value = 1.072579058299722
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.064826292770487
$ws.Range("D15").Value = 1.063060028081166
$ws.Range("E15").Value = 1.068965875959031
$ws.Range("F15").Value = 1.076896814853247
$ws.Range("I15").Value = 1.053024292026066
$ws.Range("J15").Value = 1.071158510422419
$ws.Range("K15").Value = 1.066514311687177
$ws.Range("L15").Value = 1.072399699381626
$ws.Range("M15").Value = 1.08030355048796
$ws.Range("N15").Value = 1.072679678370817
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.065623585638235
$ws.Range("D16").Value = 1.063549179177167
$ws.Range("E16").Value = 1.069639606329806
$ws.Range("F16").Value = 1.077537094412596
$ws.Range("I16").Value = 1.053241623544022
$ws.Range("J16").Value = 1.071743101185243
$ws.Range("K16").Value = 1.066891303832411
$ws.Range("L16").Value = 1.07296137144496
$ws.Range("M16").Value = 1.080832835000144
$ws.Range("N16").Value = 1.073265099319578
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.066123671241735
$ws.Range("D17").Value = 1.063856007641466
$ws.Range("E17").Value = 1.070062279154182
$ws.Range("F17").Value = 1.077938793034928
$ws.Range("I17").Value = 1.053377559162052
$ws.Range("J17").Value = 1.072109597223661
$ws.Range("K17").Value = 1.067127559492047
$ws.Range("L17").Value = 1.073313556841738
$ws.Range("M17").Value = 1.081164708340927
$ws.Range("N17").Value = 1.073632115824408
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.066415350850837
$ws.Range("D18").Value = 1.064034974400829
$ws.Range("E18").Value = 1.070308838902658
$ws.Range("F18").Value = 1.078173121855
$ws.Range("I18").Value = 1.053456707441374
$ws.Range("J18").Value = 1.07232329571907
$ws.Range("K18").Value = 1.067265283692703
$ws.Range("L18").Value = 1.073518931744298
$ws.Range("M18").Value = 1.081358236825765
$ws.Range("N18").Value = 1.073846117796192
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.066514804122674
$ws.Range("D19").Value = 1.064095997257033
$ws.Range("E19").Value = 1.070392913144472
$ws.Range("F19").Value = 1.078253026157306
$ws.Range("I19").Value = 1.053483671112242
$ws.Range("J19").Value = 1.072396149138921
$ws.Range("K19").Value = 1.067312230573122
$ws.Range("L19").Value = 1.073588951021446
$ws.Range("M19").Value = 1.081424216984298
$ws.Range("N19").Value = 1.073919074676256
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.066070018035526
$ws.Range("D20").Value = 1.063823087965199
$ws.Range("E20").Value = 1.070016928110499
$ws.Range("F20").Value = 1.077895692014474
$ws.Range("I20").Value = 1.053362989110012
$ws.Range("J20").Value = 1.072070283171679
$ws.Range("K20").Value = 1.067102219745523
$ws.Range("L20").Value = 1.073275775746736
$ws.Range("M20").Value = 1.081129106404381
$ws.Range("N20").Value = 1.073592745941963
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.064623830395117
$ws.Range("D21").Value = 1.062935820465726
$ws.Range("E21").Value = 1.068794818479654
$ws.Range("F21").Value = 1.076734254062715
$ws.Range("I21").Value = 1.052968986002794
$ws.Range("J21").Value = 1.071010006708794
$ws.Range("K21").Value = 1.066418516354981
$ws.Range("L21").Value = 1.072257035573799
$ws.Range("M21").Value = 1.080169111879454
$ws.Range("N21").Value = 1.072530963764885
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.063714430762195
$ws.Range("D22").Value = 1.062377949850026
$ws.Range("E22").Value = 1.06802661839024
$ws.Range("F22").Value = 1.07600423026704
$ws.Range("I22").Value = 1.052719990439242
$ws.Range("J22").Value = 1.070342706265932
$ws.Range("K22").Value = 1.065987924149759
$ws.Range("L22").Value = 1.071616065015377
$ws.Range("M22").Value = 1.079565090716974
$ws.Range("N22").Value = 1.071862715678854
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.064196533307432
$ws.Range("D23").Value = 1.062673688620947
$ws.Range("E23").Value = 1.068433838291242
$ws.Range("F23").Value = 1.076391209248231
$ws.Range("I23").Value = 1.052852108237117
$ws.Range("J23").Value = 1.070696518138844
$ws.Range("K23").Value = 1.066216257663589
$ws.Range("L23").Value = 1.071955898496884
$ws.Range("M23").Value = 1.079885335216089
$ws.Range("N23").Value = 1.072217030005208
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.066094261646624
$ws.Range("D24").Value = 1.063837962953465
$ws.Range("E24").Value = 1.070037420227679
$ws.Range("F24").Value = 1.077915167433871
$ws.Range("I24").Value = 1.05336957312407
$ws.Range("J24").Value = 1.072088047723434
$ws.Range("K24").Value = 1.067113669932716
$ws.Range("L24").Value = 1.073292847548438
$ws.Range("M24").Value = 1.081145193534309
$ws.Range("N24").Value = 1.07361053572142
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.068296013074606
$ws.Range("D25").Value = 1.065188995709756
$ws.Range("E25").Value = 1.071899148533524
$ws.Range("F25").Value = 1.079684611793831
$ws.Range("I25").Value = 1.053964545658562
$ws.Range("J25").Value = 1.073700003487542
$ws.Range("K25").Value = 1.068151940465834
$ws.Range("L25").Value = 1.07484238886405
$ws.Range("M25").Value = 1.082605323612142
$ws.Range("N25").Value = 1.075224780647606
